# Update crypto price table cells per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.666.17'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = "'" + '1.644.07'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'" + '215.09'
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("E6").Value = '  +1.79%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  +1.10%  '

$ws.Range("D9").Value = "'" + '0.0626'
$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").Value = "'" + '19.27'
$ws.Range("E10").Value = '  +1.44%  '

$ws.Range("D11").Value = "'" + '0.0842'
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").Value = "'" + '1.872.51'
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("E13").Value = '  +2.18%  '

$ws.Range("D14").Value = "'" + '1.599.37'
$ws.Range("E14").Value = '  -2.03%  '

$ws.Range("E15").Value = '  +2.08%  '

$ws.Range("D16").Value = "'" + '65.12'
$ws.Range("E16").Value = '  +2.67%  '

$ws.Range("D17").Value = "'" + '26.682.56'
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").Value = "'" + '0.0' + [char]0x2083 + '0744'
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("D19").Value = "'" + '216.33'
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").Value = "'" + '1.00'
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").Value = "'" + '6.27'
$ws.Range("E22").Value = '  +1.92%  '

$ws.Range("D23").Value = "'" + '9.50'
$ws.Range("E23").Value = '  +1.82%  '

$ws.Range("E24").Value = '  +15.40%  '

$ws.Range("D25").Value = "'" + '145.88'
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").Value = "'" + '7.16'
$ws.Range("E28").Value = '  +4.52%  '

$ws.Range("D29").Value = "'" + '15.72'
$ws.Range("E29").Value = '  +1.36%  '

$ws.Range("E30").Value = '  +2.45%  '

$ws.Range("D31").Value = "'" + '1.17'
$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("D33").Value = "'" + '3.04'
$ws.Range("E33").Value = '  +3.52%  '

$ws.Range("D34").Value = "'" + '1.275.50'
$ws.Range("E34").Value = '  +4.97%  '

$ws.Range("E35").Value = '  +2.59%  '

$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").Value = "'" + '0.0180'
$ws.Range("E37").Value = '  +4.52%  '

$ws.Range("D38").Value = "'" + '0.533'
$ws.Range("E38").Value = '  +7.01%  '

$ws.Range("D39").Value = "'" + '0.832'
$ws.Range("E39").Value = '  +4.96%  '

$ws.Range("D40").Value = "'" + '1.00'
$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").Value = "'" + '0.815'
$ws.Range("E41").Value = '  +2.53%  '

$ws.Range("E42").Value = '  -1.34%  '

$ws.Range("D43").Value = "'" + '5.45'
$ws.Range("E43").Value = '  +2.24%  '

$ws.Range("D44").Value = "'" + '1.782.77'
$ws.Range("E44").Value = '  +1.10%  '

$ws.Range("D45").Value = "'" + '91.31'
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").Value = "'" + '59.62'
$ws.Range("E46").Value = '  +8.62%  '

$ws.Range("D47").Value = "'" + '1.60'
$ws.Range("E47").Value = '  +2.58%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = "'" + '0.0' + [char]0x2086 + '0102'
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'" + '0.0515'
$ws.Range("E49").Value = '  +0.77%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'" + '7.81'
$ws.Range("E50").Value = '  +2.83%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = "'" + '0.0968'
$ws.Range("E51").Value = '  +2.91%  '
